$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.931.50"
$ws.Range("E2").Value = "  -1.11%  "

$ws.Range("D3").Value = "1.638.35"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'215.31"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("D9").Value = "'0.0641"
$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("D10").Value = "'19.63"
$ws.Range("E10").Value = "  -1.66%  "

$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").Value = "1.865.27"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").Value = "1.641.36"
$ws.Range("E14").Value = "  -4.86%  "

$ws.Range("E15").Value = "  -1.19%  "

$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").Value = "'62.96"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "25.937.01"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").Value = "'193.06"
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("E22").Value = "  -1.11%  "

$ws.Range("E23").Value = "  -0.91%  "

$ws.Range("D24").Value = "'144.04"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").Value = "'0.129"
$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("E32").Value = "  -1.38%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("D35").Value = "'2.44"
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("D37").Value = "1.140.30"
$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").Value = "'0.546"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("D39").Value = "'2.45"
$ws.Range("E39").Value = "  -1.38%  "

$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").Value = "'5.48"

$ws.Range("D43").Value = "'99.34"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").Value = "1.775.00"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("E46").Value = "  +2.17%  "

$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").Value = "'0.0533"
$ws.Range("E48").Value = "  +2.92%  "

$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").Value = "'7.67"
$ws.Range("E50").Value = "  -0.51%  "

$ws.Range("E51").Value = "  -0.89%  "
